# updated main GSC export data
#
# Appends the next day (2025-12-13) to the "Chart" sheet's Date/Non-HTTPS/HTTPS
# table (row 69) and leaves the "Table" sheet (Issue/Validation/Pages header)
# untouched - its values don't actually change, only their shared-string ids
# would shift in a from-scratch export.

$wb = $excel.ActiveWorkbook

$chart = $wb.Worksheets.Item("Chart")

# Figure out where the data currently ends and add the next row right after it.
$lastRow = $chart.Cells.Item($chart.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Writing the date string straight into a General-formatted cell makes Excel's
# COM layer auto-convert it to a date serial (like typing it in the UI would).
# The source data stores dates as plain text, so build the literal via a
# formula first, then paste-special it back in as a static value - that keeps
# it a plain text cell (matches how every other date cell in the sheet is
# stored) instead of a number/date cell.
$dateCell = $chart.Cells.Item($newRow, 1)
$dateCell.Formula = '="2025-12-13"'
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$chart.Cells.Item($newRow, 2).Value = 0
$chart.Cells.Item($newRow, 3).Value = 30
